$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SamplesTab" query (cell B3) is updated to a new, shorter query that
# no longer returns the Tumor status / Analyte Type columns.
$newSamplesQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001524'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value2 = $newSamplesQuery

# Update the active window's selection/scroll position to match the edited
# workbook (now showing row 3 at the top, with B3 selected).
$ws.Range("B3").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1

$wb.Save()
